$wb = $excel.ActiveWorkbook
$aboutSheet = $wb.Worksheets.Item("About")
$ws = $wb.Worksheets.Item("RTMF-freight")

# Freight HDV row: zero out shifts to rail (E) and ships (F) so that the
# remaining share is treated as freight logistics/eliminated trips (column I).
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0

# Mirror the author's final selection on the freight sheet (cell I3, which
# holds the recalculated "eliminated trips" fraction), then restore the
# originally active "About" tab.
$ws.Range("I3").Select()
$aboutSheet.Activate()
